# correção nos dados e inicio da analise PNAD 2009
#
# The original sheet had several "category header" rows (e.g. "sexo",
# "cor ou raça", "grupos de idade", "nível de instrução",
# "classes de rendimento mensal domiciliar per capita", "sem rendimento a
# menos") that carried a label in column A but no data values, plus two
# trailing footnote-only rows at the very end of the table. This edit
# removes all of those empty/footnote rows so the table becomes a clean,
# contiguous block of rows that each carry real data, and it fixes the
# pandas-generated placeholder header text in B2 ("unnamed: 1_level_1")
# to read "total" (matching B1).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the empty "section header" rows and the trailing footnote rows.
# Work from the bottom of the sheet upwards so that earlier row numbers
# stay valid while later ones are removed.
$rowsToDelete = @(36, 35, 29, 27, 19, 13, 8, 5)
foreach ($r in $rowsToDelete) {
    $ws.Rows.Item($r).Delete()
}

# Fix the mis-labelled header cell (was the pandas placeholder
# "unnamed: 1_level_1"); it should read "total" like B1.
$ws.Range("B2").Value = "total"
